$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# File path used across the new rows (same files already referenced by the
# existing "LoginPage*"/"AddEmployeePage*" rows -- note the doubled
# backslashes, matching the literal text already stored in the workbook)
$orangeHrmPath = "C:\\AutomationTesting\\OrangeHRM_Project\\src\\test\\resources\\TestData\\OrangeHRM.xlsx"
$keywordsPath  = "C:\\AutomationTesting\\OrangeHRM_Project\\src\\test\\resources\\TestData\\Keywords.xlsx"
$locatorPath   = "C:\\AutomationTesting\\OrangeHRM_Project\\src\\test\\resources\\TestData\\Locator_Data.xlsx"

# Row 8: EmployeeListPageTestData
$ws.Range("A8").Value = "EmployeeListPageTestData"
$ws.Range("B8").Value = $orangeHrmPath
$ws.Range("C8").Value = "EmployeeListPage"

# Row 9: EmployeeListPageKeywords
$ws.Range("A9").Value = "EmployeeListPageKeywords"
$ws.Range("B9").Value = $keywordsPath
$ws.Range("C9").Value = "EmployeeListPage"

# Row 10: EmployeeListPageLocators
$ws.Range("A10").Value = "EmployeeListPageLocators"
$ws.Range("B10").Value = $locatorPath
$ws.Range("C10").Value = "EmployeeListPage"

# Row 11: DashboardPageLocators
$ws.Range("A11").Value = "DashboardPageLocators"
$ws.Range("B11").Value = $locatorPath
$ws.Range("C11").Value = "DashboardPage"

# Rows 8-11 share the same row height / wrapped formatting as the other data rows
$ws.Rows("8:11").RowHeight = 28.8

# Move the active selection to B11 (matches the saved cursor position)
$ws.Range("B11").Select()
